# This script updates the "cryptos" price-tracking worksheet with a fresh
# data pull (price + 1h volume change, and a few reshuffled coin rows),
# mirroring the automated "Updated cryptos list ... with GitHub Actions"
# commit that refreshes this sheet on a schedule.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is the cell reference and its new text value. Values that look
# like plain numbers (e.g. "1.00", "9.03") are prefixed with a leading
# apostrophe so Excel stores/keeps them as text (matching the workbook's
# existing text-formatted Price column) instead of silently converting them
# to numeric values and dropping the trailing/insignificant zeros.
$updates = @(
    @{ Cell = 'D2'; Value = '70.909.44' },
    @{ Cell = 'E2'; Value = '  -0.31%  ' },
    @{ Cell = 'D3'; Value = '3.804.68' },
    @{ Cell = 'E3'; Value = '  -1.61%  ' },
    @{ Cell = 'D4'; Value = '''1.00' },
    @{ Cell = 'E4'; Value = '  -0.08%  ' },
    @{ Cell = 'D5'; Value = '''705.07' },
    @{ Cell = 'E5'; Value = '  +1.13%  ' },
    @{ Cell = 'D6'; Value = '''170.35' },
    @{ Cell = 'E6'; Value = '  -2.06%  ' },
    @{ Cell = 'D7'; Value = '3.802.89' },
    @{ Cell = 'E7'; Value = '  -1.58%  ' },
    @{ Cell = 'E8'; Value = '  +0.00%  ' },
    @{ Cell = 'E9'; Value = '  -0.80%  ' },
    @{ Cell = 'E10'; Value = '  -1.25%  ' },
    @{ Cell = 'D11'; Value = '''7.62' },
    @{ Cell = 'E11'; Value = '  +5.59%  ' },
    @{ Cell = 'E12'; Value = '  -0.46%  ' },
    @{ Cell = 'E13'; Value = '  -3.15%  ' },
    @{ Cell = 'D14'; Value = '''35.77' },
    @{ Cell = 'E14'; Value = '  -1.99%  ' },
    @{ Cell = 'D15'; Value = '4.445.72' },
    @{ Cell = 'E15'; Value = '  -1.63%  ' },
    @{ Cell = 'D16'; Value = '3.806.49' },
    @{ Cell = 'E16'; Value = '  -1.69%  ' },
    @{ Cell = 'D17'; Value = '70.881.10' },
    @{ Cell = 'E17'; Value = '  -0.45%  ' },
    @{ Cell = 'D18'; Value = '''17.38' },
    @{ Cell = 'E18'; Value = '  -1.55%  ' },
    @{ Cell = 'E19'; Value = '  -0.24%  ' },
    @{ Cell = 'E20'; Value = '  -1.73%  ' },
    @{ Cell = 'D21'; Value = '''497.21' },
    @{ Cell = 'E21'; Value = '  -0.07%  ' },
    @{ Cell = 'D22'; Value = '''10.65' },
    @{ Cell = 'E22'; Value = '  -3.90%  ' },
    @{ Cell = 'E23'; Value = '  -0.21%  ' },
    @{ Cell = 'D24'; Value = '''83.99' },
    @{ Cell = 'E24'; Value = '  -1.15%  ' },
    @{ Cell = 'D25'; Value = '''0.0000143' },
    @{ Cell = 'E25'; Value = '  -5.11%  ' },
    @{ Cell = 'D26'; Value = '3.952.57' },
    @{ Cell = 'E26'; Value = '  -1.40%  ' },
    @{ Cell = 'E27'; Value = '  -1.86%  ' },
    @{ Cell = 'D28'; Value = '''10.29' },
    @{ Cell = 'E28'; Value = '  -4.26%  ' },
    @{ Cell = 'E29'; Value = '  +0.09%  ' },
    @{ Cell = 'E30'; Value = '  -6.56%  ' },
    @{ Cell = 'E31'; Value = '  -4.57%  ' },
    @{ Cell = 'E32'; Value = '  -1.59%  ' },
    @{ Cell = 'E33'; Value = '  -3.48%  ' },
    @{ Cell = 'D34'; Value = '''29.01' },
    @{ Cell = 'E34'; Value = '  -2.34%  ' },
    @{ Cell = 'D35'; Value = '''0.174' },
    @{ Cell = 'E35'; Value = '  -4.92%  ' },
    @{ Cell = 'B36'; Value = 'Binance-PegBSC-USD' },
    @{ Cell = 'C36'; Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd' },
    @{ Cell = 'D36'; Value = '''0.995' },
    @{ Cell = 'E36'; Value = '  -0.46%  ' },
    @{ Cell = 'B37'; Value = 'RenzoRestakedETH' },
    @{ Cell = 'C37'; Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth' },
    @{ Cell = 'D37'; Value = '3.768.69' },
    @{ Cell = 'E37'; Value = '  -1.36%  ' },
    @{ Cell = 'B38'; Value = 'Aptos' },
    @{ Cell = 'C38'; Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt' },
    @{ Cell = 'D38'; Value = '''9.03' },
    @{ Cell = 'E38'; Value = '  -2.20%  ' },
    @{ Cell = 'E39'; Value = '  -3.92%  ' },
    @{ Cell = 'D40'; Value = '''2.36' },
    @{ Cell = 'E40'; Value = '  -1.42%  ' },
    @{ Cell = 'B41'; Value = 'Filecoin' },
    @{ Cell = 'C41'; Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil' },
    @{ Cell = 'D41'; Value = '''5.94' },
    @{ Cell = 'E41'; Value = '  -1.33%  ' },
    @{ Cell = 'B42'; Value = 'Mantle' },
    @{ Cell = 'C42'; Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt' },
    @{ Cell = 'D42'; Value = '''1.02' },
    @{ Cell = 'E42'; Value = '  -3.25%  ' },
    @{ Cell = 'B43'; Value = 'dogwifhat' },
    @{ Cell = 'C43'; Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif' },
    @{ Cell = 'D43'; Value = '''3.24' },
    @{ Cell = 'E43'; Value = '  -5.44%  ' },
    @{ Cell = 'B44'; Value = 'USDe' },
    @{ Cell = 'C44'; Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde' },
    @{ Cell = 'D44'; Value = '''1.00' },
    @{ Cell = 'E44'; Value = '  -0.01%  ' },
    @{ Cell = 'E45'; Value = '  +0.09%  ' },
    @{ Cell = 'D46'; Value = '''167.01' },
    @{ Cell = 'E46'; Value = '  +2.00%  ' },
    @{ Cell = 'D47'; Value = '''0.000312' },
    @{ Cell = 'E47'; Value = '  +0.37%  ' },
    @{ Cell = 'D48'; Value = '''49.05' },
    @{ Cell = 'E48'; Value = '  -0.79%  ' },
    @{ Cell = 'D49'; Value = '''416.11' },
    @{ Cell = 'E49'; Value = '  -0.11%  ' },
    @{ Cell = 'E50'; Value = '  -0.72%  ' },
    @{ Cell = 'B51'; Value = 'ONDO' },
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo' },
    @{ Cell = 'D51'; Value = '''1.35' },
    @{ Cell = 'E51'; Value = '  -2.94%  ' }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
